$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking text (e.g. "239.29") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers and
# mangles the exact textual representation (trailing zeros, scientific notation, etc).
# These cells were (and must remain) text cells in the source data.
$forceTextCells = @(
    'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D17', 'D20', 'D21', 'D22', 'D23', 'D24', 'D27', 'D28', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D46', 'D47', 'D50'
)

foreach ($cellAddr in $forceTextCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# New values for every changed cell (Coin / Link / Price / Volume(1h) columns).
$updates = @(
    @{ Cell = 'D2'; Value = '43.928.67' },
    @{ Cell = 'E2'; Value = '  -0.96%  ' },
    @{ Cell = 'D3'; Value = '2.357.28' },
    @{ Cell = 'E3'; Value = '  -0.42%  ' },
    @{ Cell = 'E4'; Value = '  +0.06%  ' },
    @{ Cell = 'D5'; Value = '0.680' },
    @{ Cell = 'E5'; Value = '  +0.75%  ' },
    @{ Cell = 'D6'; Value = '239.29' },
    @{ Cell = 'D7'; Value = '74.25' },
    @{ Cell = 'E7'; Value = '  +0.69%  ' },
    @{ Cell = 'E8'; Value = '  -0.02%  ' },
    @{ Cell = 'D9'; Value = '0.591' },
    @{ Cell = 'E9'; Value = '  +7.54%  ' },
    @{ Cell = 'D10'; Value = '0.100' },
    @{ Cell = 'E10'; Value = '  -1.68%  ' },
    @{ Cell = 'D11'; Value = '57.22' },
    @{ Cell = 'E11'; Value = '  -0.28%  ' },
    @{ Cell = 'D12'; Value = '32.02' },
    @{ Cell = 'E12'; Value = '  +6.54%  ' },
    @{ Cell = 'B13'; Value = 'TRON' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' },
    @{ Cell = 'D13'; Value = '0.108' },
    @{ Cell = 'E13'; Value = '  +0.89%  ' },
    @{ Cell = 'B14'; Value = 'Polkadot' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D14'; Value = '7.23' },
    @{ Cell = 'E14'; Value = '  +6.29%  ' },
    @{ Cell = 'D15'; Value = '2.707.43' },
    @{ Cell = 'E15'; Value = '  -0.34%  ' },
    @{ Cell = 'D16'; Value = '16.50' },
    @{ Cell = 'E16'; Value = '  -1.96%  ' },
    @{ Cell = 'D17'; Value = '0.900' },
    @{ Cell = 'E17'; Value = '  -0.45%  ' },
    @{ Cell = 'D18'; Value = '2.367.54' },
    @{ Cell = 'E18'; Value = '  -0.13%  ' },
    @{ Cell = 'D19'; Value = '43.870.92' },
    @{ Cell = 'E19'; Value = '  -1.23%  ' },
    @{ Cell = 'D20'; Value = '6.93' },
    @{ Cell = 'E20'; Value = '  +7.20%  ' },
    @{ Cell = 'D21'; Value = '0.0000101' },
    @{ Cell = 'E21'; Value = '  -1.41%  ' },
    @{ Cell = 'D22'; Value = '77.28' },
    @{ Cell = 'E22'; Value = '  -0.18%  ' },
    @{ Cell = 'D23'; Value = '258.13' },
    @{ Cell = 'E23'; Value = '  +1.44%  ' },
    @{ Cell = 'D24'; Value = '1.96' },
    @{ Cell = 'E24'; Value = '  +22.33%  ' },
    @{ Cell = 'E25'; Value = '  -0.06%  ' },
    @{ Cell = 'E26'; Value = '  -5.38%  ' },
    @{ Cell = 'D27'; Value = '2.48' },
    @{ Cell = 'E27'; Value = '  -1.27%  ' },
    @{ Cell = 'D28'; Value = '10.78' },
    @{ Cell = 'E28'; Value = '  +3.71%  ' },
    @{ Cell = 'E29'; Value = '  -0.71%  ' },
    @{ Cell = 'D30'; Value = '22.84' },
    @{ Cell = 'E30'; Value = '  +1.33%  ' },
    @{ Cell = 'D31'; Value = '175.39' },
    @{ Cell = 'E31'; Value = '  +0.81%  ' },
    @{ Cell = 'E32'; Value = '  -1.68%  ' },
    @{ Cell = 'E33'; Value = '  +2.18%  ' },
    @{ Cell = 'D34'; Value = '0.0757' },
    @{ Cell = 'E34'; Value = '  +1.91%  ' },
    @{ Cell = 'D35'; Value = '5.55' },
    @{ Cell = 'E35'; Value = '  +6.76%  ' },
    @{ Cell = 'D36'; Value = '5.20' },
    @{ Cell = 'E36'; Value = '  +0.06%  ' },
    @{ Cell = 'D37'; Value = '3.75' },
    @{ Cell = 'E37'; Value = '  -3.73%  ' },
    @{ Cell = 'D38'; Value = '6.33' },
    @{ Cell = 'E38'; Value = '  -2.16%  ' },
    @{ Cell = 'E39'; Value = '  -3.44%  ' },
    @{ Cell = 'D40'; Value = '0.0278' },
    @{ Cell = 'E40'; Value = '  +2.33%  ' },
    @{ Cell = 'D41'; Value = '0.111' },
    @{ Cell = 'E41'; Value = '  +12.48%  ' },
    @{ Cell = 'D42'; Value = '0.202' },
    @{ Cell = 'E42'; Value = '  +8.88%  ' },
    @{ Cell = 'D43'; Value = '8.97' },
    @{ Cell = 'E43'; Value = '  +1.40%  ' },
    @{ Cell = 'D44'; Value = '18.82' },
    @{ Cell = 'E44'; Value = '  -3.58%  ' },
    @{ Cell = 'E45'; Value = '  +0.03%  ' },
    @{ Cell = 'D46'; Value = '59.23' },
    @{ Cell = 'E46'; Value = '  +12.71%  ' },
    @{ Cell = 'D47'; Value = '4.73' },
    @{ Cell = 'E47'; Value = '  +5.54%  ' },
    @{ Cell = 'E48'; Value = '  +5.10%  ' },
    @{ Cell = 'E49'; Value = '  -0.66%  ' },
    @{ Cell = 'D50'; Value = '100.88' },
    @{ Cell = 'E50'; Value = '  +1.95%  ' },
    @{ Cell = 'E51'; Value = '  -0.41%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Restore the default (unstyled) cell style now that the text values are locked in,
# so the forced Text number format does not linger as visible cell formatting.
foreach ($cellAddr in $forceTextCells) {
    $ws.Range($cellAddr).Style = "Normal"
}
